$d = $word.ActiveDocument
$vt = [char]11

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# ------------------------------------------------------------------
# Body paragraph ("The intersection of science and art ...") -- edit
# from the end of the paragraph backwards so each edit only normalizes
# the runs it actually touches.
# ------------------------------------------------------------------

# Collapse the trailing three sentences (separated by manual line
# breaks) down into a single new sentence.
$oldTail = " This reciprocal relationship between science and art fosters a fertile ground for interdisciplinary exploration and discovery.$vt$vt" + `
    "In the realm of art and science convergence, artists and scientists engage in a dialogue that transcends traditional boundaries." + `
    " They share a common language of curiosity, exploration, and the pursuit of knowledge." + `
    " This dialogue enriches both disciplines, fostering a deeper understanding of the world and our place within it"
$newTail = " By actively participating in the creative process, students learn to embrace challenges, think outside the box, and find innovative solutions to problems"
Replace-Text $oldTail $newTail

Replace-Text " Artists' unique perspectives and methodologies prompt scientists to approach problems from novel angles, leading to breakthroughs and innovations" " They learn to appreciate different perspectives and understand the nuances of communication"

Replace-Text " The creative process inherent in art encourages scientists to think outside the conventional boundaries of their disciplines" " In analyzing and interpreting works of art, students develop the ability to think critically and engage in meaningful discourse"

Replace-Text "Conversely, science also finds a muse in art" "Furthermore, the arts foster essential critical thinking and problem-solving skills"

Replace-Text " These natural marvels become the raw materials of artistic creation, transformed into paintings, sculptures, and installations that capture the essence of scientific concepts" " The arts provide a safe haven for exploration, where students can experiment with different forms of expression, allowing them to develop their own artistic voices"

Replace-Text " Artists draw inspiration from the colors of the aurora borealis, the patterns of snowflakes, and the intricate geometry of DNA" " Through engagement with visual arts, music, theater, and dance, students can cultivate their unique perspectives and find creative outlets for self-expression"

Replace-Text "In the tapestry of this convergence, art finds its muse in the wonders of the natural world and the intricacies of scientific phenomena" "In a world increasingly dominated by technology and quantitative reasoning, the arts offer a refuge for the imaginative mind"

Replace-Text " This essay delves into the fascinating landscape of this convergence, exploring how science inspires art and how art, in turn, informs science" " This essay explores the invaluable role of arts education in shaping well-rounded, creative, and expressive individuals"

Replace-Text " It is a space where artistic expression and scientific inquiry dance in harmonious unison, each enriching the other" " As educators, it is our duty to cultivate a comprehensive educational experience that values the arts alongside the traditional subjects"

Replace-Text "The intersection of science and art is a realm where boundaries blur, and creativity and knowledge converge" "The arts, in their myriad forms, have long been recognized for their ability to enrich our lives and expand our horizons"

# ------------------------------------------------------------------
# Summary paragraph -- also edited back to front.
# ------------------------------------------------------------------

$oldSummaryTail = " Simultaneously, art informs science, challenging conventional thinking and stimulating innovative approaches to problem-solving." + `
    " The reciprocal relationship between these disciplines fosters interdisciplinary exploration and discovery, leading to a deeper understanding of the world and our place within it"
$newSummaryTail = " By recognizing the importance of the arts in education, we empower students to become well-rounded individuals who are equipped to make meaningful contributions to society"
Replace-Text $oldSummaryTail $newSummaryTail

Replace-Text " Science inspires art, providing a wealth of ideas and wonders that fuel artistic expression" " It cultivates creativity, critical thinking, and empathy, while fostering a sense of personal and cultural identity"

Replace-Text "The convergence of science and art is a dynamic and ever-evolving realm where creativity and knowledge intertwine" "In essence, arts education provides students with the tools and skills necessary to navigate the complexities of an ever-changing world"

# ------------------------------------------------------------------
# Title / byline / email
# ------------------------------------------------------------------

Replace-Text "xyz" "com"
Replace-Text "james@emailcentral" "albright@schoolmail"
Replace-Text "ethan" "jessica"
Replace-Text "Ethan James" "Jessica Albright"
Replace-Text "Unveiling the Convergence of Science and Art" "Beyond the Classroom: The Value of Arts Education"

# ------------------------------------------------------------------
# Append a trailing empty paragraph at the very end of the document.
# ------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$endOfDoc = $lastPara.Range.Text
$ok = $d.Content.Find.Execute("By recognizing the importance of the arts in education, we empower students to become well-rounded individuals who are equipped to make meaningful contributions to society.", $true, $false, $false, $false, $false, $true, 1, $false, "By recognizing the importance of the arts in education, we empower students to become well-rounded individuals who are equipped to make meaningful contributions to society.^p", 2)
Write-Output "trailing paragraph insert ok=$ok"
